$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(344, 44418, 0, 1, 16.63616702711695),
    @(345, 44419, 1, 2, 33.27233405423391),
    @(346, 44420, 0, 1, 16.63616702711695),
    @(347, 44421, 1, 2, 33.27233405423391),
    @(348, 44422, 0, 2, 33.27233405423391),
    @(349, 44423, 0, 2, 33.27233405423391),
    @(350, 44424, 0, 2, 33.27233405423391),
    @(351, 44425, 1, 3, 49.90850108135086),
    @(352, 44426, 0, 2, 33.27233405423391),
    @(353, 44427, 1, 3, 49.90850108135086),
    @(354, 44428, 3, 5, 83.18083513558476),
    @(355, 44429, 0, 5, 83.18083513558476),
    @(356, 44430, 1, 6, 99.81700216270171),
    @(357, 44431, 0, 6, 99.81700216270171)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    # Apply same style as column A on row 343 (date format with border)
    $ws.Range("A343").Copy()
    $ws.Range("A$r").PasteSpecial(-4122) # xlPasteFormats
}
